# Insert a new weekly price record for "Haba" (Macroferia Regional de Talca)
# as row 12, pushing the existing rows 12-38 down to 13-39.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..38 down by one (this also grows the used range / dimension
# to A1:R39 and carries the row-12 number-format style, s="2", for column D).
$ws.Rows.Item(12).Insert()

# Populate the newly-opened row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value  = 5
$ws.Cells.Item(12, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value  = "Maule"
$ws.Cells.Item(12, 4).Value  = 44487
$ws.Cells.Item(12, 5).Value  = 7
$ws.Cells.Item(12, 6).Value  = 100112026
$ws.Cells.Item(12, 7).Value  = "Haba"
$ws.Cells.Item(12, 8).Value  = "Sin especificar"
$ws.Cells.Item(12, 9).Value  = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 8000
$ws.Cells.Item(12, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 16).Value = 320
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
